# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The sheet's "K" column (G) is rebuilt from the refreshed data source. Row 35
# (r=34 in the sheet, i.e. A-index 33) is unchanged; every other data row
# (2-52) gets a new strikeout (K) count written in place of the stale value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 2
    6  = 1
    7  = 2
    8  = 2
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 4
    22 = 4
    23 = 0
    24 = 1
    25 = 2
    26 = 4
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 1
    34 = 3
    36 = 3
    37 = 3
    38 = 1
    39 = 1
    40 = 3
    41 = 3
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 2
    48 = 2
    49 = 1
    50 = 2
    51 = 2
    52 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

$excel.Calculate()
